# Update the "取得日時" (retrieved-at) timestamp in column A for rows 2-11
# on the "ランサーズ" sheet from "2025-10-14 06:28:19" to "2025-10-14 06:33:46".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("ランサーズ")

$oldValue = "2025-10-14 06:28:19"
$newValue = "2025-10-14 06:33:46"

for ($row = 2; $row -le 11; $row++) {
    $cell = $ws.Cells.Item($row, 1)
    if ($cell.Text -eq $oldValue) {
        $cell.Value = $newValue
    }
}
